$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — copy H1's formatting (bold,
# bordered, centered) so the new header cells share style index with the
# existing header row instead of minting a fresh style.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I ("I0") and J ("IF"), rows 2-11.
$iValues = @(8, 7, 7, 8, 8, 7, 8, 6, 6, 3)
$jValues = @(9, 8, 8, 8, 8, 7, 8, 6, 6, 4)

for ($r = 0; $r -lt 10; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
